$d = $word.ActiveDocument

# Change 1: VM hardware/OS description in the Evaluation section
$d.Content.Find.Execute(
    "three Microsoft Azure B2S VMs with 2 vCPUs, ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "three Microsoft Azure Standard DS1 VMs running Windows 2019 Datacenter Server with 1 vCPUs, ",
    2
)

# Change 2: Add "East" after Japan in the list of regions
$d.Content.Find.Execute(
    "situated in Central Canada, Southern UK and Japan.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "situated in Central Canada, Southern UK and Japan East.",
    2
)
